$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new work log entry in row 20
$ws.Range("B20").Value = "Tackling the outliers problems + continuing with regression analysis"
$ws.Range("C20").Value = 1

# Reflect the final selection left by the author after entering the data
$ws.Range("E22").Select()
